$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 14: A empty, B="nowa7" (same text as row10/13), C=38, D="adam", E="krzywo"
$ws.Range("A14").Font.Bold = $false
$ws.Range("B14").Value = "nowa7"
$ws.Range("C14").Value = 38
$ws.Range("D14").Value = "adam"
$ws.Range("E14").Value = "krzywo"

# Row 15: A empty, B="nowa8", C=39, D="adam", E="krzywo"
$ws.Range("A15").Font.Bold = $false
$ws.Range("B15").Value = "nowa8"
$ws.Range("C15").Value = 39
$ws.Range("D15").Value = "adam"
$ws.Range("E15").Value = "krzywo"

# Row 16: A empty, B="nowa9", C=40, D="adam", E="krzywo"
$ws.Range("A16").Font.Bold = $false
$ws.Range("B16").Value = "nowa9"
$ws.Range("C16").Value = 40
$ws.Range("D16").Value = "adam"
$ws.Range("E16").Value = "krzywo"

# Row 17: A empty, B="nowa10", C=41, D="adam", E="krzywo"
$ws.Range("A17").Font.Bold = $false
$ws.Range("B17").Value = "nowa10"
$ws.Range("C17").Value = 41
$ws.Range("D17").Value = "adam"
$ws.Range("E17").Value = "krzywo"

# Row 18: A empty, B="nowa11", C=42, D="adam", E="krzywo"
$ws.Range("A18").Font.Bold = $false
$ws.Range("B18").Value = "nowa11"
$ws.Range("C18").Value = 42
$ws.Range("D18").Value = "adam"
$ws.Range("E18").Value = "krzywo"

# Row 19: A=16, B="nowa12", C=43, D="adam", E="krzywo"
$ws.Range("A19").Value = 16
$ws.Range("B19").Value = "nowa12"
$ws.Range("C19").Value = 43
$ws.Range("D19").Value = "adam"
$ws.Range("E19").Value = "krzywo"

# Row 20: A empty, B="nowa13", C=44, D="adam", E="krzywo"
$ws.Range("A20").Font.Bold = $false
$ws.Range("B20").Value = "nowa13"
$ws.Range("C20").Value = 44
$ws.Range("D20").Value = "adam"
$ws.Range("E20").Value = "krzywo"
